# Applies the "Paie Info Repaprtition Bug Found" correction:
# - Month label: Juilet -> Aout
# - Code: 4159 -> 4190 (both occurrences, in both table rows)
# - Gross amount: 74 862 000,00 -> 75 420 000,00 (both occurrences)
# - 5% amount: 3 743 100,00 -> 3 771 000,00 (both occurrences)
# - Amount in words: updated to match the new total

$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

$d.Content.Find.Execute("Juilet", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Aout", $wdReplaceAll)

$d.Content.Find.Execute("4159", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4190", $wdReplaceAll)

$d.Content.Find.Execute("74 862 000,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "75 420 000,00", $wdReplaceAll)

$d.Content.Find.Execute("3 743 100,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3 771 000,00", $wdReplaceAll)

$d.Content.Find.Execute("TROIS MILLIONS SEPT CENT QUARANTE-TROIS MILLE CENT ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TROIS MILLIONS SEPT CENT SOIXANTE ET ONZE MILLE  ", $wdReplaceAll)
